$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 25; existing rows 25-38 shift down to 26-39.
$ws.Rows(25).Insert()

# Populate the newly inserted row 25 with the new weekly record.
$ws.Cells.Item(25, 1).Value = 9
$ws.Cells.Item(25, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(25, 3).Value = "Metropolitana"
$ws.Cells.Item(25, 4).Value = 44960
$ws.Cells.Item(25, 5).Value = 13
$ws.Cells.Item(25, 6).Value = 100112010
$ws.Cells.Item(25, 7).Value = "Achicoria"
$ws.Cells.Item(25, 8).Value = "Sin especificar"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 70
$ws.Cells.Item(25, 11).Value = 7000
$ws.Cells.Item(25, 12).Value = 7000
$ws.Cells.Item(25, 13).Value = 7000
$ws.Cells.Item(25, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(25, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(25, 16).Value = 438
$ws.Cells.Item(25, 17).Value = 16
$ws.Cells.Item(25, 18).Value = "Hortaliza"
